# Update NATMI LR-pair metrics with recomputed TPM-based values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1.0
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.002991666666666667
$ws.Range("H2").Value = 0.008975
$ws.Range("I2").Value = 0.0003566413595017623
$ws.Range("J2").Value = 0.0003566413595017623
$ws.Range("M2").Value = 8.676671
$ws.Range("N2").Value = 26.030013
$ws.Range("O2").Value = 0.1325240072999665
$ws.Range("P2").Value = 0.1325240072999665
$ws.Range("Q2").Value = 0.02595770740833334
$ws.Range("R2").Value = 0.233619366675
$ws.Range("S2").Value = 0.00004726354213008153
$ws.Range("T2").Value = 0.00004726354213008151
$ws.Range("E3").Value = 1.0
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.002991666666666667
$ws.Range("H3").Value = 0.008975
$ws.Range("I3").Value = 0.0003566413595017623
$ws.Range("J3").Value = 0.0003566413595017623
$ws.Range("M3").Value = 37.74750533333334
$ws.Range("O3").Value = 0.5765403197090441
$ws.Range("P3").Value = 0.576540319709044
$ws.Range("Q3").Value = 0.1129279534555556
$ws.Range("R3").Value = 1.0163515811
$ws.Range("S3").Value = 0.0002056181234286142
$ws.Range("T3").Value = 0.0002056181234286141
$ws.Range("E4").Value = 1.0
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.002991666666666667
$ws.Range("H4").Value = 0.008975
$ws.Range("I4").Value = 0.0003566413595017623
$ws.Range("J4").Value = 0.0003566413595017623
$ws.Range("M4").Value = 19.04827033333333
$ws.Range("N4").Value = 57.144811
$ws.Range("O4").Value = 0.2909356729909895
$ws.Range("P4").Value = 0.2909356729909895
$ws.Range("Q4").Value = 0.05698607541388889
$ws.Range("R4").Value = 0.512874678725
$ws.Range("S4").Value = 0.0001037596939430666
$ws.Range("T4").Value = 0.0001037596939430666
$ws.Range("I5").Value = 0.9971069332391614
$ws.Range("J5").Value = 0.9971069332391616
$ws.Range("M5").Value = 8.676671
$ws.Range("N5").Value = 26.030013
$ws.Range("O5").Value = 0.1325240072999665
$ws.Range("P5").Value = 0.1325240072999665
$ws.Range("Q5").Value = 72.57321490699066
$ws.Range("R5").Value = 653.158934162916
$ws.Range("S5").Value = 0.1321406064994338
$ws.Range("T5").Value = 0.1321406064994338
$ws.Range("I6").Value = 0.9971069332391614
$ws.Range("J6").Value = 0.9971069332391616
$ws.Range("M6").Value = 37.74750533333334
$ws.Range("O6").Value = 0.5765403197090441
$ws.Range("P6").Value = 0.576540319709044
$ws.Range("S6").Value = 0.5748723500738105
$ws.Range("T6").Value = 0.5748723500738105
$ws.Range("I7").Value = 0.9971069332391614
$ws.Range("J7").Value = 0.9971069332391616
$ws.Range("M7").Value = 19.04827033333333
$ws.Range("N7").Value = 57.144811
$ws.Range("O7").Value = 0.2909356729909895
$ws.Range("P7").Value = 0.2909356729909895
$ws.Range("Q7").Value = 159.3231109612724
$ws.Range("R7").Value = 1433.907998651452
$ws.Range("S7").Value = 0.2900939766659171
$ws.Range("T7").Value = 0.2900939766659171
$ws.Range("G8").Value = 0.02127666666666667
$ws.Range("H8").Value = 0.06383
$ws.Range("I8").Value = 0.002536425401336767
$ws.Range("J8").Value = 0.002536425401336767
$ws.Range("M8").Value = 8.676671
$ws.Range("N8").Value = 26.030013
$ws.Range("O8").Value = 0.1325240072999665
$ws.Range("P8").Value = 0.1325240072999665
$ws.Range("Q8").Value = 0.1846106366433334
$ws.Range("R8").Value = 1.66149572979
$ws.Range("S8").Value = 0.0003361372584025742
$ws.Range("T8").Value = 0.0003361372584025742
$ws.Range("G9").Value = 0.02127666666666667
$ws.Range("H9").Value = 0.06383
$ws.Range("I9").Value = 0.002536425401336767
$ws.Range("J9").Value = 0.002536425401336767
$ws.Range("M9").Value = 37.74750533333334
$ws.Range("O9").Value = 0.5765403197090441
$ws.Range("P9").Value = 0.576540319709044
$ws.Range("Q9").Value = 0.8031410884755555
$ws.Range("R9").Value = 7.22826979628
$ws.Range("S9").Value = 0.00146235151180484
$ws.Range("T9").Value = 0.00146235151180484
$ws.Range("G10").Value = 0.02127666666666667
$ws.Range("H10").Value = 0.06383
$ws.Range("I10").Value = 0.002536425401336767
$ws.Range("J10").Value = 0.002536425401336767
$ws.Range("M10").Value = 19.04827033333333
$ws.Range("N10").Value = 57.144811
$ws.Range("O10").Value = 0.2909356729909895
$ws.Range("P10").Value = 0.2909356729909895
$ws.Range("Q10").Value = 0.4052836984588889
$ws.Range("R10").Value = 3.64755328613
$ws.Range("S10").Value = 0.000737936631129353
$ws.Range("T10").Value = 0.000737936631129353
